$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completeness")

$ws.Range("A3").Value = "water temperature"
$ws.Range("A5").Value = "DO concentration"
$ws.Range("A6").Value = "sp conductivity"
$ws.Range("A9").Value = "orthoP"
$ws.Range("A13").Value = "chlorophyll a"

$ws.Range("A13").Select()
